# Added working test for fund ratio computation
# - Rename "Demo Fund 1" to "Demo Fund 2" on the CapitalCall sheet
# - Update the active selection on the CapitalCall sheet to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Update the fund name values (shared string "Demo Fund 1" -> "Demo Fund 2")
$ws.Range("A2").Value = "Demo Fund 2"
$ws.Range("A3").Value = "Demo Fund 2"

# Move/reset the selection on the CapitalCall sheet to A2
$ws.Activate()
$ws.Range("A2").Select()
